$d = $word.ActiveDocument

$d.Content.Find.Execute(
    "rendelkezik felhasználói fiókkal",
    $true, $false, $false, $false, $false,
    $true, 1, $false, "be van jelentkezve", 2
)
